$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.037425168578215
$ws.Cells.Item(2, 4).Value = 1.038039982128679
$ws.Cells.Item(2, 5).Value = 1.036169267099454
$ws.Cells.Item(2, 6).Value = 1.044811884571918
$ws.Cells.Item(2, 9).Value = 1.034999749369511
$ws.Cells.Item(2, 10).Value = 1.042528110566738
$ws.Cells.Item(2, 11).Value = 1.040829149743927
$ws.Cells.Item(2, 12).Value = 1.038963777913633
$ws.Cells.Item(2, 13).Value = 1.047581883668407
$ws.Cells.Item(2, 14).Value = 1.044008620063444
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.038866572148391
$ws.Cells.Item(3, 4).Value = 1.039318777423955
$ws.Cells.Item(3, 5).Value = 1.037410084754151
$ws.Cells.Item(3, 6).Value = 1.046360826011885
$ws.Cells.Item(3, 9).Value = 1.035298800602293
$ws.Cells.Item(3, 10).Value = 1.043611417103615
$ws.Cells.Item(3, 11).Value = 1.041916637253422
$ws.Cells.Item(3, 12).Value = 1.040012997134534
$ws.Cells.Item(3, 13).Value = 1.048940216838673
$ws.Cells.Item(3, 14).Value = 1.045093465019804
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.039797804498898
$ws.Cells.Item(4, 4).Value = 1.040145137839241
$ws.Cells.Item(4, 5).Value = 1.038211961001082
$ws.Cells.Item(4, 6).Value = 1.047361802769458
$ws.Cells.Item(4, 9).Value = 1.035490247004678
$ws.Cells.Item(4, 10).Value = 1.044310544086983
$ws.Cells.Item(4, 11).Value = 1.042618664907793
$ws.Cells.Item(4, 12).Value = 1.040690344716441
$ws.Cells.Item(4, 13).Value = 1.049817370166303
$ws.Cells.Item(4, 14).Value = 1.04579358484368
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.040188953862373
$ws.Cells.Item(5, 4).Value = 1.040492280396769
$ws.Cells.Item(5, 5).Value = 1.03854883163704
$ws.Cells.Item(5, 6).Value = 1.047782311666149
$ws.Cells.Item(5, 9).Value = 1.035570239641436
$ws.Cells.Item(5, 10).Value = 1.04460402043635
$ws.Cells.Item(5, 11).Value = 1.042913407331592
$ws.Cells.Item(5, 12).Value = 1.040974731036776
$ws.Cells.Item(5, 13).Value = 1.050185706278103
$ws.Cells.Item(5, 14).Value = 1.046087477963126
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.04025460978249
$ws.Cells.Item(6, 4).Value = 1.04055055210501
$ws.Cells.Item(6, 5).Value = 1.038605379880258
$ws.Cells.Item(6, 6).Value = 1.04785289948536
$ws.Cells.Item(6, 9).Value = 1.035583641987834
$ws.Cells.Item(6, 10).Value = 1.044653270947318
$ws.Cells.Item(6, 11).Value = 1.042962873153002
$ws.Cells.Item(6, 12).Value = 1.041022459153215
$ws.Cells.Item(6, 13).Value = 1.050247527148862
$ws.Cells.Item(6, 14).Value = 1.046136798415469
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.039803032386051
$ws.Cells.Item(7, 4).Value = 1.040149777387605
$ws.Cells.Item(7, 5).Value = 1.038216463212509
$ws.Cells.Item(7, 6).Value = 1.047367422806171
$ws.Cells.Item(7, 9).Value = 1.035491317799131
$ws.Cells.Item(7, 10).Value = 1.044314467241821
$ws.Cells.Item(7, 11).Value = 1.042622604796897
$ws.Cells.Item(7, 12).Value = 1.040694146150865
$ws.Cells.Item(7, 13).Value = 1.049822293530065
$ws.Cells.Item(7, 14).Value = 1.045797513569849
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.037912602520934
$ws.Cells.Item(8, 4).Value = 1.038472389138461
$ws.Cells.Item(8, 5).Value = 1.036588820591203
$ws.Cells.Item(8, 6).Value = 1.045335628647058
$ws.Cells.Item(8, 9).Value = 1.035101242629446
$ws.Cells.Item(8, 10).Value = 1.04289460408273
$ws.Cells.Item(8, 11).Value = 1.041197015726116
$ws.Cells.Item(8, 12).Value = 1.039318693805622
$ws.Cells.Item(8, 13).Value = 1.048041311171852
$ws.Cells.Item(8, 14).Value = 1.044375634042211
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.034569990348785
$ws.Cells.Item(9, 4).Value = 1.035507908658906
$ws.Cells.Item(9, 5).Value = 1.033712700306588
$ws.Cells.Item(9, 6).Value = 1.041745113345173
$ws.Cells.Item(9, 9).Value = 1.03439802950165
$ws.Cells.Item(9, 10).Value = 1.04037826935242
$ws.Cells.Item(9, 11).Value = 1.03867209253788
$ws.Cells.Item(9, 12).Value = 1.036882757648072
$ws.Cells.Item(9, 13).Value = 1.04488906302147
$ws.Cells.Item(9, 14).Value = 1.041855725827956
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.032333457372417
$ws.Cells.Item(10, 4).Value = 1.033525395288165
$ws.Cells.Item(10, 5).Value = 1.031789584597355
$ws.Cells.Item(10, 6).Value = 1.039344078145551
$ws.Cells.Item(10, 9).Value = 1.033918456324619
$ws.Cells.Item(10, 10).Value = 1.038690743960822
$ws.Cells.Item(10, 11).Value = 1.036979858695315
$ws.Cells.Item(10, 12).Value = 1.035250294248527
$ws.Cells.Item(10, 13).Value = 1.04277777111728
$ws.Cells.Item(10, 14).Value = 1.040165803956741
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.031362988717605
$ws.Cells.Item(11, 4).Value = 1.032665400187495
$ws.Cells.Item(11, 5).Value = 1.030955429960757
$ws.Cells.Item(11, 6).Value = 1.038302551236731
$ws.Cells.Item(11, 9).Value = 1.03370821836924
$ws.Cells.Item(11, 10).Value = 1.037957590915937
$ws.Cells.Item(11, 11).Value = 1.036244909840083
$ws.Cells.Item(11, 12).Value = 1.034541335875967
$ws.Cells.Item(11, 13).Value = 1.041861139787614
$ws.Cells.Item(11, 14).Value = 1.039431609750437
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.031002199160671
$ws.Cells.Item(12, 4).Value = 1.032345719646597
$ws.Cells.Item(12, 5).Value = 1.030645366242753
$ws.Cells.Item(12, 6).Value = 1.037915392610999
$ws.Cells.Item(12, 9).Value = 1.03362973691553
$ws.Cells.Item(12, 10).Value = 1.037684891888914
$ws.Cells.Item(12, 11).Value = 1.0359715804265
$ws.Cells.Item(12, 12).Value = 1.03427767723711
$ws.Cells.Item(12, 13).Value = 1.041520288973432
$ws.Cells.Item(12, 14).Value = 1.039158523459517
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.031079604105317
$ws.Cells.Item(13, 4).Value = 1.032414303178621
$ws.Cells.Item(13, 5).Value = 1.030711886089527
$ws.Cells.Item(13, 6).Value = 1.037998452685385
$ws.Cells.Item(13, 9).Value = 1.033646589124132
$ws.Cells.Item(13, 10).Value = 1.037743403751297
$ws.Cells.Item(13, 11).Value = 1.036030225847659
$ws.Cells.Item(13, 12).Value = 1.034334247481298
$ws.Cells.Item(13, 13).Value = 1.04159341964952
$ws.Cells.Item(13, 14).Value = 1.039217118415456
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.031333172201275
$ws.Cells.Item(14, 4).Value = 1.032638980230198
$ws.Cells.Item(14, 5).Value = 1.030929804553052
$ws.Cells.Item(14, 6).Value = 1.03827055453668
$ws.Cells.Item(14, 9).Value = 1.033701739030451
$ws.Cells.Item(14, 10).Value = 1.03793505718169
$ws.Cells.Item(14, 11).Value = 1.036222323249202
$ws.Cells.Item(14, 12).Value = 1.03451954833084
$ws.Cells.Item(14, 13).Value = 1.041832972613137
$ws.Cells.Item(14, 14).Value = 1.039409044015703
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.03148936197128
$ws.Cells.Item(15, 4).Value = 1.032777379165906
$ws.Cells.Item(15, 5).Value = 1.031064041765092
$ws.Cells.Item(15, 6).Value = 1.038438166872401
$ws.Cells.Item(15, 9).Value = 1.033735667001777
$ws.Cells.Item(15, 10).Value = 1.038053091526242
$ws.Cells.Item(15, 11).Value = 1.036340636001066
$ws.Cells.Item(15, 12).Value = 1.034633675731641
$ws.Cells.Item(15, 13).Value = 1.041980519406759
$ws.Cells.Item(15, 14).Value = 1.039527245982563
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.03239782057628
$ws.Cells.Item(16, 4).Value = 1.033582437029965
$ws.Cells.Item(16, 5).Value = 1.031844913970439
$ws.Cells.Item(16, 6).Value = 1.039413160817484
$ws.Cells.Item(16, 9).Value = 1.033932354592111
$ws.Cells.Item(16, 10).Value = 1.03873934890185
$ws.Cells.Item(16, 11).Value = 1.03702858795065
$ws.Cells.Item(16, 12).Value = 1.035297300933701
$ws.Cells.Item(16, 13).Value = 1.042838553053729
$ws.Cells.Item(16, 14).Value = 1.040214477922359
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.032967121794028
$ws.Cells.Item(17, 4).Value = 1.034087007911223
$ws.Cells.Item(17, 5).Value = 1.032334346446341
$ws.Cells.Item(17, 6).Value = 1.04002424349354
$ws.Cells.Item(17, 9).Value = 1.034055039368979
$ws.Cells.Item(17, 10).Value = 1.039169161457003
$ws.Cells.Item(17, 11).Value = 1.037459528627136
$ws.Cells.Item(17, 12).Value = 1.035713011610117
$ws.Cells.Item(17, 13).Value = 1.043376118974868
$ws.Cells.Item(17, 14).Value = 1.040644900860641
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.033298989783721
$ws.Cells.Item(18, 4).Value = 1.034381166109306
$ws.Cells.Item(18, 5).Value = 1.03261968630772
$ws.Cells.Item(18, 6).Value = 1.040380498576159
$ws.Cells.Item(18, 9).Value = 1.03412635057645
$ws.Cells.Item(18, 10).Value = 1.039419628519783
$ws.Cells.Item(18, 11).Value = 1.037710677137703
$ws.Cells.Item(18, 12).Value = 1.035955287262571
$ws.Cells.Item(18, 13).Value = 1.043689438382383
$ws.Cells.Item(18, 14).Value = 1.040895723615379
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.03341211524143
$ws.Cells.Item(19, 4).Value = 1.034481441301784
$ws.Cells.Item(19, 5).Value = 1.03271695663852
$ws.Cells.Item(19, 6).Value = 1.040501942253472
$ws.Cells.Item(19, 9).Value = 1.034150623728403
$ws.Cells.Item(19, 10).Value = 1.039504991587637
$ws.Cells.Item(19, 11).Value = 1.037796276531812
$ws.Cells.Item(19, 12).Value = 1.036037863020545
$ws.Cells.Item(19, 13).Value = 1.043796232815186
$ws.Cells.Item(19, 14).Value = 1.040981207908581
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.03290606147336
$ws.Cells.Item(20, 4).Value = 1.034032887738404
$ws.Cells.Item(20, 5).Value = 1.032281849274647
$ws.Cells.Item(20, 6).Value = 1.039958698698474
$ws.Cells.Item(20, 9).Value = 1.034041902198258
$ws.Cells.Item(20, 10).Value = 1.039123071023996
$ws.Cells.Item(20, 11).Value = 1.037413314732428
$ws.Cells.Item(20, 12).Value = 1.035668430657549
$ws.Cells.Item(20, 13).Value = 1.043318467472106
$ws.Cells.Item(20, 14).Value = 1.040598744973933
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.031258511461805
$ws.Cells.Item(21, 4).Value = 1.032572825119114
$ws.Cells.Item(21, 5).Value = 1.03086563916342
$ws.Cells.Item(21, 6).Value = 1.03819043538057
$ws.Cells.Item(21, 9).Value = 1.033685509533135
$ws.Cells.Item(21, 10).Value = 1.037878630379784
$ws.Cells.Item(21, 11).Value = 1.036165764695706
$ws.Cells.Item(21, 12).Value = 1.034464990727619
$ws.Cells.Item(21, 13).Value = 1.041762440603487
$ws.Cells.Item(21, 14).Value = 1.039352537081266
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.03022080872731
$ws.Cells.Item(22, 4).Value = 1.031653433505894
$ws.Cells.Item(22, 5).Value = 1.029973926717891
$ws.Cells.Item(22, 6).Value = 1.037076980418695
$ws.Cells.Item(22, 9).Value = 1.033459175457147
$ws.Cells.Item(22, 10).Value = 1.03709403797858
$ws.Cells.Item(22, 11).Value = 1.03537942946051
$ws.Cells.Item(22, 12).Value = 1.033706486342302
$ws.Cells.Item(22, 13).Value = 1.040781941783138
$ws.Cells.Item(22, 14).Value = 1.038566830468859
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.030771090212987
$ws.Cells.Item(23, 4).Value = 1.032140954599057
$ws.Cells.Item(23, 5).Value = 1.030446764262891
$ws.Cells.Item(23, 6).Value = 1.037667406100838
$ws.Cells.Item(23, 9).Value = 1.033579373988866
$ws.Cells.Item(23, 10).Value = 1.037510172398047
$ws.Cells.Item(23, 11).Value = 1.035796467627348
$ws.Cells.Item(23, 12).Value = 1.034108761545693
$ws.Cells.Item(23, 13).Value = 1.041301930586564
$ws.Cells.Item(23, 14).Value = 1.038983555846933
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.032933652609144
$ws.Cells.Item(24, 4).Value = 1.034057342777025
$ws.Cells.Item(24, 5).Value = 1.032305570912959
$ws.Cells.Item(24, 6).Value = 1.039988316123679
$ws.Cells.Item(24, 9).Value = 1.034047839088935
$ws.Cells.Item(24, 10).Value = 1.039143898034111
$ws.Cells.Item(24, 11).Value = 1.037434197457168
$ws.Cells.Item(24, 12).Value = 1.035688575493444
$ws.Cells.Item(24, 13).Value = 1.043344518426963
$ws.Cells.Item(24, 14).Value = 1.040619601560791
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.035435535491759
$ws.Cells.Item(25, 4).Value = 1.036275364529681
$ws.Cells.Item(25, 5).Value = 1.034457227819428
$ws.Cells.Item(25, 6).Value = 1.04267460937057
$ws.Cells.Item(25, 9).Value = 1.03458171633203
$ws.Cells.Item(25, 10).Value = 1.041030535831557
$ws.Cells.Item(25, 11).Value = 1.039326400626716
$ws.Cells.Item(25, 12).Value = 1.037513983042459
$ws.Cells.Item(25, 13).Value = 1.045705689293819
$ws.Cells.Item(25, 14).Value = 1.042508918600309
